$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 390-394 (U column "xp" tweaks, AA/AB/AC "gold" range tweaks) ---
$ws.Range("U390").Value = 49
$ws.Range("U391").Value = 54
$ws.Range("U392").Value = 55
$ws.Range("U393").Value = 58
$ws.Range("AA393:AC393").Value = 260000000000000
$ws.Range("U394").Value = 63
$ws.Range("AA394:AC394").Value = 267000000000000

# --- Append six new monster rows (395-400) ---
$newMonsters = @(
    @{ Row = 395; Name = "Rebirthed Abomination";     Stat8 = 274000000000000; Stat4 = 17.8333333333333;    DamageStat = "dex";   Xp = 66; Gold = 530000000; HealthRange = "360476637656600-385476476445500"; AE = 8.4999999999999893; AF = 11.85; AG = 8.6999999999999993; AH = 8.8999999999999897; AI = 15.9 },
    @{ Row = 396; Name = "The Creators Left Hand";    Stat8 = 281000000000000; Stat4 = 18.233333333333299;   DamageStat = "focus"; Xp = 68; Gold = 540000000; HealthRange = "420476637656600-455476476445500"; AE = 8.5599999999999898; AF = 11.94; AG = 8.76;              AH = 8.9599999999999902; AI = 16.02 },
    @{ Row = 397; Name = "The Creators Right Hand";   Stat8 = 288000000000000; Stat4 = 18.633333333333301;   DamageStat = "int";   Xp = 71; Gold = 550000000; HealthRange = "450476637656600-485476476445500"; AE = 8.6199999999999903; AF = 12.03; AG = 8.82;              AH = 9.0199999999999907; AI = 16.14 },
    @{ Row = 398; Name = "The Celestial Creator";     Stat8 = 295000000000000; Stat4 = 19.033333333333299;   DamageStat = "chr";   Xp = 73; Gold = 560000000; HealthRange = "530476637656600-545476476445500"; AE = 8.6799999999999908; AF = 12.12; AG = 8.8800000000000008;  AH = 9.0799999999999894; AI = 16.260000000000002 },
    @{ Row = 399; Name = "Cosmic Satan";              Stat8 = 302000000000000; Stat4 = 19.433333333333302;   DamageStat = "chr";   Xp = 75; Gold = 570000000; HealthRange = "560476637656600-585476476445500"; AE = 8.7399999999999896; AF = 12.21; AG = 8.94;              AH = 9.1399999999999899; AI = 16.38 },
    @{ Row = 400; Name = "Heretical Prophet of Time"; Stat8 = 309000000000000; Stat4 = 19.8333333333333;     DamageStat = "dex";   Xp = 78; Gold = 580000000; HealthRange = "640476637656600-685476476445500"; AE = 8.7999999999999901; AF = 12.3;  AG = 9;                  AH = 9.1999999999999904; AI = 16.5 }
)

foreach ($m in $newMonsters) {
    $r = $m.Row

    $ws.Range("A$r").Value = $m.Name

    $ws.Range("B$r`:I$r").Value = $m.Stat8

    $ws.Range("J$r`:M$r").Value = $m.Stat4

    $ws.Range("Q$r").Value = 1
    $ws.Range("R$r").Value = 1
    $ws.Range("S$r").Value = 99999
    $ws.Range("T$r").Value = $m.DamageStat
    $ws.Range("U$r").Value = $m.Xp
    $ws.Range("V$r").Value = 0.001
    $ws.Range("W$r").Value = $m.Gold
    $ws.Range("X$r").Value = 0
    $ws.Range("Y$r").Value = $m.HealthRange
    $ws.Range("Z$r").Value = $m.HealthRange

    $ws.Range("AA$r`:AC$r").Value = $m.Stat8

    $ws.Range("AD$r").Value = 1
    $ws.Range("AE$r").Value = $m.AE
    $ws.Range("AF$r").Value = $m.AF
    $ws.Range("AG$r").Value = $m.AG
    $ws.Range("AH$r").Value = $m.AH
    $ws.Range("AI$r").Value = $m.AI
    $ws.Range("AJ$r").Value = 1

    $ws.Range("AM$r").Value = "Purgatory"
}

# --- Match the final UI selection state from the authored edit ---
$ws.Range("A401").Select()
